# Apply updates to the NBA MVP PPG dataframe:
#  - Rename column A header "year" -> "season_ending_year"
#  - Add new column K "calendar_year" holding the year as a number
#    (same values as column A, but numeric instead of text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header of column A
$ws.Range("A1").Value = "season_ending_year"

# Add the new header for column K and copy the header formatting/style
# from column A (bold font, border, centered alignment)
$ws.Range("K1").Value = "calendar_year"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Fill column K with the numeric calendar year, mirroring column A's
# text year values for each data row
for ($r = 2; $r -le 11; $r++) {
    $yearText = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 11).Value = [int]$yearText
}
